$wb = $excel.ActiveWorkbook

# --- 展览 (sheet1) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 383
$ws.Range("F5").Value = 24
$ws.Range("F6").Value = 1246
$ws.Range("F7").Value = 450
$ws.Range("F8").Value = 102
$ws.Range("F9").Value = 191
$ws.Range("F10").Value = 153
$ws.Range("F11").Value = 175
$ws.Range("F12").Value = 1053
$ws.Range("F13").Value = 5
$ws.Range("F14").Value = 270
$ws.Range("F15").Value = 188
$ws.Range("F16").Value = 1504
$ws.Range("F17").Value = 549
$ws.Range("F18").Value = 231
$ws.Range("F19").Value = 349
$ws.Range("F20").Value = 113
$ws.Range("F21").Value = 822
$ws.Range("F22").Value = 1154
$ws.Range("F24").Value = 1915
$ws.Range("F25").Value = 2662
$ws.Range("F26").Value = 1443
$ws.Range("F27").Value = 65
$ws.Range("F28").Value = 39
$ws.Range("F29").Value = 414
$ws.Range("F30").Value = 419
$ws.Range("F31").Value = 1245
$ws.Range("F32").Value = 825
$ws.Range("F33").Value = 1362
$ws.Range("F36").Value = 786
$ws.Range("F37").Value = 615
$ws.Range("F38").Value = 674
$ws.Range("F39").Value = 847
$ws.Range("F40").Value = 365
$ws.Range("F41").Value = 249
# --- 演出 (sheet2) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("I9").Value = "//i2.hdslb.com/bfs/openplatform/202405/xfdusgJP1715147982566.jpeg"
$ws.Range("F15").Value = 632
# --- 全部类型 (sheet4) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value = 383
$ws.Range("F8").Value = 25
$ws.Range("F11").Value = 1246
$ws.Range("F12").Value = 450
$ws.Range("F13").Value = 102
$ws.Range("F14").Value = 191
$ws.Range("F16").Value = 153
$ws.Range("F17").Value = 175
$ws.Range("F18").Value = 1053
$ws.Range("F19").Value = 270
$ws.Range("F20").Value = 188
$ws.Range("F21").Value = 1504
$ws.Range("F22").Value = 550
$ws.Range("F23").Value = 231
$ws.Range("F24").Value = 349
$ws.Range("F28").Value = 1154
$ws.Range("F29").Value = 2662
$ws.Range("F30").Value = 1443
$ws.Range("F31").Value = 65
$ws.Range("F34").Value = 414
$ws.Range("F35").Value = 419
$ws.Range("F36").Value = 1245
$ws.Range("F39").Value = 825
$ws.Range("F40").Value = 1362
$ws.Range("F41").Value = 786
$ws.Range("F42").Value = 615
$ws.Range("F43").Value = 674
$ws.Range("F44").Value = 847
$ws.Range("F45").Value = 365
$ws.Range("F48").Value = 249

